$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.442.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.27"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4646"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3721"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8875"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07931"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.95%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.860.79"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.411"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.603"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.45"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008900"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.89"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.474.13"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.150"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.063.10"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.94%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.33"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.51"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.079"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.143"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.75"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08896"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.029"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7555"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.15%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.489"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.660"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +10.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.081"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01965"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05260"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.990"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.150"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5191"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.348"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4860"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.32"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.88"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.654"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06253"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.77"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.73%  "
